$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COCKTAIL & BAR")

# Insert 6 blank rows above the old row 95 (soft-drink list gets duplicated
# here), which pushes the existing rows 95-104 down to 101-110.
$ws.Rows("95:100").Insert() | Out-Null

# Re-use the existing "analcolici" block (rows 12-17: iced tea peach/lemon,
# pomegranate, orange soda, lemon soda, chinotto) as the content for the
# newly inserted rows, copying both values and formatting.
$ws.Range("A12:AW17").Copy() | Out-Null
$ws.Range("A95").PasteSpecial(-4122) | Out-Null
$ws.Range("A12:AW17").Copy() | Out-Null
$ws.Range("A95").PasteSpecial(-4163) | Out-Null

# Make "COCKTAIL & BAR" the active sheet/tab (was "BIRRA & SIDRO" before).
$ws.Activate() | Out-Null
$ws.Range("E96").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 72
$excel.ActiveWindow.ScrollColumn = 1
